$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.960.70"
$ws.Range("E2").Value = "  -0.49%  "
$ws.Range("D3").Value = "2.329.70"
$ws.Range("E3").Value = "  +1.03%  "
$ws.Range("E4").Value = "  +0.35%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "302.98"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.31%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "95.85"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.93%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.502"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.66%  "
$ws.Range("E8").Value = "  +0.36%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.497"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.33%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "34.19"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.01%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "19.14"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.21%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0785"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.61%  "
$ws.Range("E13").Value = "  +2.87%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.73"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.58%  "
$ws.Range("D15").Value = "2.696.92"
$ws.Range("E15").Value = "  +1.19%  "
$ws.Range("D16").Value = "2.306.26"
$ws.Range("E16").Value = "  +0.26%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.791"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.67%  "
$ws.Range("D18").Value = "42.908.52"
$ws.Range("E18").Value = "  -0.20%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.11"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.79%  "
$ws.Range("E20").Value = "  +1.59%  "
$ws.Range("E21").Value = "  -0.91%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "68.05"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.20%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "236.68"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.41%  "
$ws.Range("E24").Value = "  +2.81%  "
$ws.Range("E25").Value = "  +0.06%  "
$ws.Range("E26").Value = "  -0.51%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "24.64"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.49%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.30"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.92%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.13"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.33%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "31.61"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.72%  "
$ws.Range("E31").Value = "  +0.29%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "139.89"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -15.60%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.00"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.08%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "17.71"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.13%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0701"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.11%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.39"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.86%  "
$ws.Range("B37").Value = "WEMIXToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.30"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.94%  "
$ws.Range("B38").Value = "ARBITRUM"
$ws.Range("C38").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.80"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.99%  "
$ws.Range("E39").Value = "  -0.18%  "
$ws.Range("B40").Value = "EnergySwap"
$ws.Range("C40").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "22.41"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +22.88%  "
$ws.Range("B41").Value = "LidoDAOToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.74"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.51%  "
$ws.Range("E42").Value = "  -0.98%  "
$ws.Range("D43").Value = "1.932.10"
$ws.Range("E43").Value = "  -3.44%  "
$ws.Range("E44").Value = "  -0.46%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.03"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -4.44%  "
$ws.Range("E46").Value = "  -2.53%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.74"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.60%  "
$ws.Range("B48").Value = "RocketPoolETH"
$ws.Range("C48").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D48").Value = "2.562.51"
$ws.Range("E48").Value = "  +1.23%  "
$ws.Range("B49").Value = "HuobiToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.87"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.53%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "53.54"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.32%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "72.78"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.07%  "
